$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 6 (hunk 1)
$ws.Range("H6").Value = 91
$ws.Range("I6").Value = 86.5
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 259.5
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -147.5
$ws.Range("N6").Value = -524

# row 11 (hunk 2)
$ws.Range("H11").Value = 128.5
$ws.Range("I11").Value = 128.5
$ws.Range("K11").Value = 128.5
$ws.Range("M11").Value = 11.5

# row 80 (hunk 3)
$ws.Range("H80").Value = 643.4761999999999
$ws.Range("I80").Value = 664.4
$ws.Range("J80").Value = 636.9375
$ws.Range("K80").Value = 1993.2
$ws.Range("L80").Value = 1910.8125
$ws.Range("M80").Value = -995.1999999999998
$ws.Range("N80").Value = -3906.8125

# row 83 (hunk 4)
$ws.Range("H83").Value = 643.4761999999999
$ws.Range("I83").Value = 664.4
$ws.Range("J83").Value = 636.9375
$ws.Range("K83").Value = 5979.599999999999
$ws.Range("L83").Value = 5732.4375
$ws.Range("M83").Value = -987.5999999999995
$ws.Range("N83").Value = -15716.4375

# row 98 (hunk 5)
$ws.Range("H98").Value = 703.8889
$ws.Range("I98").Value = 644.6667
$ws.Range("K98").Value = 644.6667
$ws.Range("M98").Value = 853.3333

# row 100 (hunk 6)
$ws.Range("H100").Value = 916153.4
$ws.Range("I100").Value = 2502171.8
$ws.Range("J100").Value = 9857.143
$ws.Range("K100").Value = 2502171.8
$ws.Range("L100").Value = 9857.143
$ws.Range("M100").Value = -2501630.8
$ws.Range("N100").Value = -10939.143

# row 122 (hunk 7)
$ws.Range("H122").Value = 703.8889
$ws.Range("I122").Value = 644.6667
$ws.Range("K122").Value = 1934.0001
$ws.Range("M122").Value = 515.9999

# row 132 (hunk 8)
$ws.Range("H132").Value = 875.0540999999999
$ws.Range("I132").Value = 875.0540999999999
$ws.Range("K132").Value = 2625.1623
$ws.Range("M132").Value = -95.16229999999996

# row 137 (hunk 9)
$ws.Range("H137").Value = 2215.7222
$ws.Range("I137").Value = 1267.75
$ws.Range("K137").Value = 3803.25
$ws.Range("M137").Value = -1253.25

# row 138 (hunk 10)
$ws.Range("H138").Value = 4144.3794
$ws.Range("I138").Value = 3499.3333
$ws.Range("J138").Value = 4312.6523
$ws.Range("K138").Value = 10497.9999
$ws.Range("L138").Value = 12937.9569
$ws.Range("M138").Value = -5357.999899999999
$ws.Range("N138").Value = -23217.9569

$ws = $wb.Worksheets.Item("ARM")
# row 2 (hunk 11)
$ws.Range("H2").Value = 2329.3333
$ws.Range("J2").Value = 2249
$ws.Range("L2").Value = 2249
$ws.Range("N2").Value = -2475

# row 61 (hunk 12)
$ws.Range("H61").Value = 2294.1
$ws.Range("I61").Value = 1686
$ws.Range("K61").Value = 1686
$ws.Range("M61").Value = -1474

# row 116 (hunk 13)
$ws.Range("H116").Value = 2329.3333
$ws.Range("J116").Value = 2249
$ws.Range("L116").Value = 2249
$ws.Range("N116").Value = -6837

# row 136 (hunk 14)
$ws.Range("H136").Value = 2294.1
$ws.Range("I136").Value = 1686
$ws.Range("K136").Value = 5058
$ws.Range("M136").Value = -2508

$ws = $wb.Worksheets.Item("BSM")
# row 3 (hunk 15)
$ws.Range("H3").Value = 2329.3333
$ws.Range("J3").Value = 2249
$ws.Range("L3").Value = 2249
$ws.Range("N3").Value = -2477

# row 20 (hunk 16)
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# row 134 (hunk 17)
$ws.Range("H134").Value = 2883.2
$ws.Range("I134").Value = 2885.25
$ws.Range("K134").Value = 8655.75
$ws.Range("M134").Value = -6120.75

$ws = $wb.Worksheets.Item("CRP")
# row 20 (hunk 18)
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# row 22 (hunk 19)
$ws.Range("H22").Value = 731.1667
$ws.Range("I22").Value = 741.4
$ws.Range("K22").Value = 741.4
$ws.Range("M22").Value = -391.4

# row 30 (hunk 20)
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# row 31 (hunk 21)
$ws.Range("H31").Value = 3388.65
$ws.Range("I31").Value = 2405.8333
$ws.Range("J31").Value = 4862.875
$ws.Range("K31").Value = 2405.8333
$ws.Range("L31").Value = 4862.875
$ws.Range("M31").Value = -2110.8333
$ws.Range("N31").Value = -5452.875

# row 34 (hunk 22)
$ws.Range("H34").Value = 3388.65
$ws.Range("I34").Value = 2405.8333
$ws.Range("J34").Value = 4862.875
$ws.Range("K34").Value = 2405.8333
$ws.Range("L34").Value = 4862.875
$ws.Range("M34").Value = -2203.8333
$ws.Range("N34").Value = -5266.875

# row 58 (hunk 23)
$ws.Range("H58").Value = 2598.7896
$ws.Range("I58").Value = 2251.2
$ws.Range("J58").Value = 2985
$ws.Range("K58").Value = 2251.2
$ws.Range("L58").Value = 2985
$ws.Range("M58").Value = -2048.2
$ws.Range("N58").Value = -3391

# row 122 (hunk 24)
$ws.Range("H122").Value = 2799.75
$ws.Range("I122").Value = 599.5
$ws.Range("K122").Value = 1798.5
$ws.Range("M122").Value = 651.5

# row 124 (hunk 25)
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# row 128 (hunk 26)
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# row 131 (hunk 27)
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# row 132 (hunk 28)
$ws.Range("H132").Value = 4213.5557
$ws.Range("I132").Value = 4230.273
$ws.Range("J132").Value = 4187.2856
$ws.Range("K132").Value = 12690.819
$ws.Range("L132").Value = 12561.8568
$ws.Range("M132").Value = -10160.819
$ws.Range("N132").Value = -17621.8568

# row 133 (hunk 29)
$ws.Range("H133").Value = 124000
$ws.Range("J133").Value = 124000
$ws.Range("L133").Value = 124000
$ws.Range("N133").Value = -129060

# row 134 (hunk 30)
$ws.Range("H134").Value = 3604
$ws.Range("I134").Value = 3604
$ws.Range("K134").Value = 10812
$ws.Range("M134").Value = -8277

# row 136 (hunk 31)
$ws.Range("H136").Value = 2598.7896
$ws.Range("I136").Value = 2251.2
$ws.Range("J136").Value = 2985
$ws.Range("K136").Value = 6753.599999999999
$ws.Range("L136").Value = 8955
$ws.Range("M136").Value = -4203.599999999999
$ws.Range("N136").Value = -14055

# row 137 (hunk 32)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# row 138 (hunk 33)
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

# row 139 (hunk 34)
$ws.Range("H139").Value = 118900
$ws.Range("J139").Value = 118900
$ws.Range("L139").Value = 118900
$ws.Range("N139").Value = -129180

# row 140 (hunk 35)
$ws.Range("H140").Value = 124900
$ws.Range("J140").Value = 124900
$ws.Range("L140").Value = 124900
$ws.Range("N140").Value = -135260

# row 141 (hunk 36)
$ws.Range("H141").Value = 52353.25
$ws.Range("J141").Value = 52353.25
$ws.Range("L141").Value = 52353.25
$ws.Range("N141").Value = -62713.25

$ws = $wb.Worksheets.Item("CUL")
# row 6 (hunk 37)
$ws.Range("H6").Value = 349
$ws.Range("I6").Value = 349
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1047
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -934
$ws.Range("N6").ClearContents()

# row 140 (hunk 38)
$ws.Range("H140").Value = 977.55554
$ws.Range("I140").Value = 977.55554
$ws.Range("K140").Value = 2932.66662
$ws.Range("M140").Value = 2247.33338

# row 141 (hunk 39)
$ws.Range("H141").Value = 9664.333000000001
$ws.Range("I141").Value = 9664.333000000001
$ws.Range("K141").Value = 28992.999
$ws.Range("M141").Value = -23812.999

$ws = $wb.Worksheets.Item("GSM")
# row 43 (hunk 40)
$ws.Range("H43").Value = 22612.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 22612.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 22612.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -22914.5

# row 132 (hunk 41)
$ws.Range("H132").Value = 3612.6667
$ws.Range("I132").Value = 3282.5
$ws.Range("J132").Value = 3832.7778
$ws.Range("K132").Value = 9847.5
$ws.Range("L132").Value = 11498.3334
$ws.Range("M132").Value = -7317.5
$ws.Range("N132").Value = -16558.3334

$ws = $wb.Worksheets.Item("LTW")
# row 16 (hunk 42)
$ws.Range("H16").Value = 414.77777
$ws.Range("I16").Value = 428.2857
$ws.Range("J16").Value = 367.5
$ws.Range("K16").Value = 428.2857
$ws.Range("L16").Value = 367.5
$ws.Range("M16").Value = -258.2857
$ws.Range("N16").Value = -707.5

# row 76 (hunk 43)
$ws.Range("H76").Value = 40288
$ws.Range("J76").Value = 40288
$ws.Range("L76").Value = 40288
$ws.Range("N76").Value = -40964

# row 79 (hunk 44)
$ws.Range("H79").Value = 40288
$ws.Range("J79").Value = 40288
$ws.Range("L79").Value = 40288
$ws.Range("N79").Value = -42628

# row 132 (hunk 45)
$ws.Range("H132").Value = 5277.1113
$ws.Range("I132").Value = 5070.857
$ws.Range("K132").Value = 15212.571
$ws.Range("M132").Value = -12682.571

$ws = $wb.Worksheets.Item("WVR")
# row 68 (hunk 46)
$ws.Range("H68").Value = 10320
$ws.Range("J68").Value = 10320
$ws.Range("L68").Value = 10320
$ws.Range("N68").Value = -11942

# row 71 (hunk 47)
$ws.Range("H71").Value = 10320
$ws.Range("J71").Value = 10320
$ws.Range("L71").Value = 30960
$ws.Range("N71").Value = -39072

# row 81 (hunk 48)
$ws.Range("H81").Value = 5979.8335
$ws.Range("I81").Value = 3376.4
$ws.Range("K81").Value = 6752.8
$ws.Range("M81").Value = -5691.8

# row 84 (hunk 49)
$ws.Range("H84").Value = 5979.8335
$ws.Range("I84").Value = 3376.4
$ws.Range("K84").Value = 33764
$ws.Range("M84").Value = -28460

# row 132 (hunk 50)
$ws.Range("H132").Value = 3279.5454
$ws.Range("I132").Value = 2466.7693
$ws.Range("K132").Value = 7400.3079
$ws.Range("M132").Value = -4870.3079
